# Build site at 2022-09-26 16:07:08 UTC
# Applies the restructuring of the LOM3229 "Programa/Syllabus/Avaliacao" block:
#  - Row 10 (Objetivos / B & C) now holds the "Durval Rodrigues Junior" docente line.
#  - Rows 13-25 are reshuffled: a couple of long paragraphs (syllabus text and the
#    bibliography list) are dropped, and every row below shifts up, producing the
#    (slightly jumbled) label/value pairing that is literally present in the target
#    OOXML. Rows 26-27 disappear entirely, shrinking the sheet to A1:C25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Drop the two trailing rows so the used range becomes A1:C25.
#    (delete the same row index twice: once the first Delete() fires, what used
#    to be row 27 slides into row 26's place)
# ---------------------------------------------------------------------------
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(26).Delete()

# ---------------------------------------------------------------------------
# 2) Row 10 - Objetivos: value changes to the "Durval" docente line.
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = '6495737 - Durval Rodrigues Junior'
$ws.Range("C10").Value = '6495737 - Durval Rodrigues Junior'

# ---------------------------------------------------------------------------
# 3) Row 13 - becomes "Programa resumido:" with a stray "01/01/2012" value
#    (exactly what the target workbook contains).
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = "01/01/2012"
$ws.Range("C13").Value = "01/01/2012"
$ws.Rows.Item(13).RowHeight = 60

# ---------------------------------------------------------------------------
# 4) Row 14 - becomes "Short syllabus:" with no B/C values any more.
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Rows.Item(14).RowHeight = 60

# ---------------------------------------------------------------------------
# 5) Row 15 - becomes "Programa:" paired with the "Durval" docente line again.
# ---------------------------------------------------------------------------
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '6495737 - Durval Rodrigues Junior'
$ws.Range("C15").Value = '6495737 - Durval Rodrigues Junior'
$ws.Rows.Item(15).RowHeight = 120

# ---------------------------------------------------------------------------
# 6) Row 16 - becomes "Syllabus:", no B/C values.
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").ClearContents()
$ws.Range("C16").ClearContents()
$ws.Rows.Item(16).RowHeight = 120

# ---------------------------------------------------------------------------
# 7) Row 17 - becomes "Avaliação:" with no B/C, and reverts to the default
#    (non-custom) row height.
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Rows.Item(17).AutoFit()

# ---------------------------------------------------------------------------
# 8) Row 18 - becomes "Método:" paired with the "Paulo Atsushi Suzuki" line.
# ---------------------------------------------------------------------------
$ws.Range("A18").Value = 'Método:'
$ws.Range("B18").Value = '1643715 - Paulo Atsushi Suzuki'
$ws.Range("C18").Value = '1643715 - Paulo Atsushi Suzuki'
$ws.Rows.Item(18).RowHeight = 60

# ---------------------------------------------------------------------------
# 9) Row 19 - becomes "Critério:" paired with the "Experimentos..." text.
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B19").Value = 'Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo.'
$ws.Range("C19").Value = 'Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo.'
$ws.Rows.Item(19).RowHeight = 60

# ---------------------------------------------------------------------------
# 10) Row 20 - becomes "Norma de recuperação:" paired with the "Média
#     aritmética..." grading-average text.
# ---------------------------------------------------------------------------
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3'
$ws.Range("C20").Value = 'Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3'
$ws.Rows.Item(20).RowHeight = 60

# ---------------------------------------------------------------------------
# 11) Row 21 - becomes "Bibliografia:" paired with the "Aplicação de uma
#     prova..." makeup-exam text, and grows to a 120pt row.
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C21").Value = 'Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Rows.Item(21).RowHeight = 120

# ---------------------------------------------------------------------------
# 12) Row 22 - becomes "Requisitos:" only, default row height.
# ---------------------------------------------------------------------------
$ws.Range("A22").Value = 'Requisitos:'
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows.Item(22).AutoFit()

# ---------------------------------------------------------------------------
# 13) Row 23 - the "Requisito" list moves up: LOB1021 line, no A value, 30pt row.
# ---------------------------------------------------------------------------
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOB1021 -  Física IV  (Requisito)`n"
$ws.Range("C23").Value = "LOB1021 -  Física IV  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

# ---------------------------------------------------------------------------
# 14) Row 24 - LOM3016 requisito line, 30pt row.
# ---------------------------------------------------------------------------
$ws.Range("A24").ClearContents()
$ws.Range("B24").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
$ws.Range("C24").Value = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"
$ws.Rows.Item(24).RowHeight = 30

# ---------------------------------------------------------------------------
# 15) Row 25 - LOM3246 "Indicação de Conjunto" line, 30pt row.
# ---------------------------------------------------------------------------
$ws.Range("B25").Value = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)`n"
$ws.Range("C25").Value = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)`n"
$ws.Rows.Item(25).RowHeight = 30

Write-Output "LOM3229 sheet restructuring applied"
